$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matchup rows for spring 23 week 15 (appended after existing row 1875)
$data = @(
    @(4, 1, 6, 2),
    @(5, 2, 6, 0),
    @(5, 2, 5, 0),
    @(3, 1, 2, 2),
    @(5, 2, 6, 0),
    @(3, 0, 3, 3),
    @(4, 0, 4, 3),
    @(5, 0, 6, 3),
    @(5, 2, 4, 1),
    @(4, 0, 3, 3),
    @(3, 2, 4, 1),
    @(6, 3, 5, 0),
    @(3, 1, 4, 2),
    @(7, 0, 5, 3),
    @(4, 2, 5, 0),
    @(5, 0, 4, 2),
    @(3, 3, 3, 0),
    @(2, 1, 4, 2),
    @(7, 3, 5, 0),
    @(5, 0, 3, 2),
    @(3, 1, 3, 2),
    @(5, 0, 5, 3),
    @(6, 1, 7, 2),
    @(4, 0, 5, 2)
)

$startRow = 1876
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}

$lastRow = $startRow + $data.Length - 1
$nextRow = $lastRow + 1

$ws.Range("A" + $nextRow).Select()
